$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# This workbook lists, for each worker, the "Periodo Mora" rows.
# Originally rows 16-124 hold worker OLMEDO DE JESUS CASTRO CAMPEON
# with periods descending from 2507 to 1607 (109 rows), followed by
# three rows for RAFAEL ENRIQUE JIMENEZ GONGORA (rows 125-127) and one
# row for CARLOS ALBERTO ARZUZA DIAZ (row 128).
#
# The update replaces that trailing block: it removes the RAFAEL and
# CARLOS rows, re-sorts OLMEDO's period rows in ascending order, and
# appends one new period (2508) for OLMEDO so the data now spans
# 1607-2508 (110 rows total) for a single worker.
# -----------------------------------------------------------------

# Row 128 (CARLOS's row) carries the special "closing" border style used
# for the last data row of the table. Stash a copy of that formatting in
# a scratch row far below the data (out of the way of the row deletion
# below) so it can be re-applied to the new final row once the RAFAEL /
# CARLOS rows are removed.
$ws.Range("B128:J128").Copy() | Out-Null
$ws.Range("B500:J500").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Drop the three RAFAEL rows (126:128) and the one CARLOS row - row 125
# (the first RAFAEL row) is left in place and overwritten below; this
# nets the table down by three rows overall, sliding the blank spacer
# rows and the signature footer (previously 133/134) up to 130/131.
# The scratch row used above (500) shifts up to 497 accordingly.
$ws.Rows("126:128").Delete()

# Re-apply the stashed "closing" formatting onto the new last data row.
$ws.Range("B497:J497").Copy() | Out-Null
$ws.Range("B125:J125").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B497:J497").Clear() | Out-Null

# Re-populate the worker data rows (16-125) for OLMEDO, periods 1607..2508
# ascending, keeping the same Salario Basico (43480) and Valor Mora (1087000).
$year = 16
$month = 7
for ($r = 16; $r -le 125; $r++) {
    $period = ("{0:D2}{1:D2}" -f $year, $month)
    $ws.Cells.Item($r, 2).Value2 = "CC"
    $ws.Cells.Item($r, 3).Value2 = "10236900"
    $ws.Cells.Item($r, 4).Value2 = "OLMEDO DE JESUS CASTRO CAMPEON"
    $ws.Cells.Item($r, 5).Value2 = $period
    $ws.Cells.Item($r, 6).Value2 = 43480
    $ws.Cells.Item($r, 7).Value2 = 1087000

    $month = $month + 1
    if ($month -eq 13) {
        $month = 1
        $year = $year + 1
    }
}

# Header summary: one worker now, 110 periods of mora (1607-2508).
$ws.Range("C13:D13").Value2 = 1
$ws.Range("F13:J13").Value2 = 110
$ws.Range("E11:J11").Value2 = 4782800

$wb.Save()
